$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Duplicate column D (values + formatting) into a new column E. This keeps
#    each existing row's original "Attendance" formatting alive in E before
#    we repurpose D's own formatting below.
# ---------------------------------------------------------------------------
$ws.Range("D1:D15").Copy($ws.Range("E1:E15"))

# ---------------------------------------------------------------------------
# 2) Re-label the headers: D becomes the short "Att" column, E becomes the
#    full "Attendance" column.
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = "Att"
$ws.Range("E1").Value = "Attendance"

# ---------------------------------------------------------------------------
# 3) Column D reverts to the plain/default look (same as columns A-C) now
#    that its old look lives on in column E.
# ---------------------------------------------------------------------------
$ws.Range("A2:A15").Copy()
$ws.Range("D2:D15").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4) Recompute column E's values from column D: Present -> TRUE, everything
#    else (Absent / null / undefined) mirrors the column D text.
# ---------------------------------------------------------------------------
$ws.Range("E2").Value = $true
$ws.Range("E3").Value = $true
$ws.Range("E4").Value = $true
$ws.Range("E5").Value = "Absent"
$ws.Range("E6").Value = "Absent"

# row 7 was sparsely populated (only Name/Age) - rebuild it fully as a
# "Present" record, using row 2 (already normalized) as the formatting
# template, then fill in the real values.
$ws.Range("A2:E2").Copy($ws.Range("A7:E7"))
$ws.Range("A7").Value = "Parker"
$ws.Range("B7").Value = 10
$ws.Range("C7").Value = 17
$ws.Range("D7").Value = "Present"
$ws.Range("E7").Value = $true

$ws.Range("E8").Value = "null"

$ws.Range("A9").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("C9").Value = "null"
$ws.Range("E9").Value = "undefined"

$ws.Range("B10").Value = "olive"
$ws.Range("E10").Value = $true

$ws.Range("E11").Value = $true
$ws.Range("E12").Value = $true
$ws.Range("E13").Value = "Absent"
$ws.Range("E14").Value = "Absent"
$ws.Range("E15").Value = "null"

# ---------------------------------------------------------------------------
# 5) New attendance rows 16-20, built from existing rows with the matching
#    Present/Absent formatting so the new cells inherit correct styling.
# ---------------------------------------------------------------------------
$ws.Range("A2:E2").Copy($ws.Range("A16:E16"))
$ws.Range("A16").Value = "Tim"
$ws.Range("B16").Value = 10
$ws.Range("C16").Value = 16
$ws.Range("D16").Value = "Present"
$ws.Range("E16").Value = $true

$ws.Range("A2:E2").Copy($ws.Range("A17:E17"))
$ws.Range("A17").Value = "Parker"
$ws.Range("B17").Value = 10
$ws.Range("C17").Value = 17
$ws.Range("D17").Value = "Present"
$ws.Range("E17").Value = $true

$ws.Range("A2:E2").Copy($ws.Range("A18:E18"))
$ws.Range("A18").Value = "Stocks"
$ws.Range("B18").Value = 9
$ws.Range("C18").Value = 13
$ws.Range("D18").Value = "Present"
$ws.Range("E18").Value = $true

$ws.Range("A5:E5").Copy($ws.Range("A19:E19"))
$ws.Range("A19").Value = "Trucker"
$ws.Range("B19").Value = 9
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = "Absent"
$ws.Range("E19").Value = "Absent"

$ws.Range("A5:E5").Copy($ws.Range("A20:E20"))
$ws.Range("A20").Value = "Hennesy"
$ws.Range("B20").Value = 9
$ws.Range("C20").Value = 18
$ws.Range("D20").Value = "Absent"
$ws.Range("E20").Value = "Absent"

# ---------------------------------------------------------------------------
# 6) New trailing row 21 - plain-formatted (like the null/undefined rows),
#    no Name value, just Class/Age/Att/Attendance.
# ---------------------------------------------------------------------------
$ws.Range("A14:D14").Copy($ws.Range("B21:E21"))
$ws.Range("B21").Value = "Tim Rose"
$ws.Range("C21").Value = 10
$ws.Range("D21").Value = 4
$ws.Range("E21").Value = 4
